$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.680.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.096.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'342.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "'0.5176"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.4369"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.09%  "
$ws.Range("D9").Value = "'53.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "'0.09210"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "'24.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.79%  "
$ws.Range("D13").Value = "'6.763"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'2.041.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.02%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'8.157"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "'102.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("D17").Value = "'0.00001152"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'1.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'21.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'0.06672"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'6.203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "'29.703.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'12.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").Value = "'2.273.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "'161.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'2.484"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").Value = "'133.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").Value = "'1.125"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.46%  "
$ws.Range("D32").Value = "'1.678"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("D33").Value = "'0.1051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").Value = "'6.188"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("D35").Value = "'3.951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("D36").Value = "'6.354"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.93%  "
$ws.Range("D37").Value = "'10.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("D38").Value = "'0.02572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'0.6975"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'12.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").Value = "'0.2213"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "'0.6775"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.73%  "
$ws.Range("D45").Value = "'14.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'2.320"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "'3.616"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "'1.201"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").Value = "'1.214"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'81.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.60%  "
